$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 22:20"

# Row 20: Tenerife - update Casos totales, Recuperados, Muertes
$ws.Range("B20").Value = 1490
$ws.Range("D20").Value = 1323
$ws.Range("E20").Value = 73

# Rows 21 and 22: Salamanca/Asturias swap places (labels and data)
$ws.Range("A21").Value = "Asturias"
$ws.Range("B21").Value = 1433
$ws.Range("C21").Value = 154
$ws.Range("D21").Value = 1209
$ws.Range("E21").Value = 70

$ws.Range("A22").Value = "Salamanca"
$ws.Range("B22").Value = 1413
$ws.Range("C22").Value = 272
$ws.Range("D22").Value = 986
$ws.Range("E22").Value = 155

# Row 44: Gran Canaria - update Recuperados
$ws.Range("D44").Value = 1323

# Rows 53-56: Melilla/La Palma/Lanzarote/Ceuta rotate places (labels and data)
$ws.Range("A53").Value = "Melilla"
$ws.Range("B53").Value = 70
$ws.Range("C53").Value = 4
$ws.Range("D53").Value = 64
$ws.Range("E53").Value = 2

$ws.Range("A54").Value = "La Palma"
$ws.Range("B54").Value = 69
$ws.Range("C54").Value = 94
$ws.Range("D54").Value = 1323
$ws.Range("E54").Value = 2

$ws.Range("A55").Value = "Lanzarote"
$ws.Range("B55").Value = 65
$ws.Range("C55").Value = 94
$ws.Range("D55").Value = 1323
$ws.Range("E55").Value = 3

$ws.Range("A56").Value = "Ceuta"
$ws.Range("B56").Value = 62
$ws.Range("C56").Value = 0
$ws.Range("D56").Value = 61
$ws.Range("E56").Value = 1

# Rows 58, 61, 63: update Recuperados
$ws.Range("D58").Value = 1323
$ws.Range("D61").Value = 1323
$ws.Range("D63").Value = 1323
